# Auto-generated-assisted edit script
# Applies cell-value updates to the "leve profit" data sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR)
# per the authoritative diff. Values are plain numbers (t="n"), not formulas, so we just
# poke .Value on each target cell; a couple of cells are removed/added entirely, handled via
# ClearContents() / a fresh .Value assignment respectively.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 185.66667
$ws.Range("I33").Value = 202.71428
$ws.Range("J33").Value = 126
$ws.Range("K33").Value = 202.71428
$ws.Range("L33").Value = 126
$ws.Range("M33").Value = 26.28572
$ws.Range("N33").Value = -584

# Row 121: Mindful Medicine / Tincture of Mind
$ws.Range("H121").Value = 595
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# Row 127: Liquid Competence / Competent Craftsman's Draught
$ws.Range("H127").Value = 909.8125
$ws.Range("I127").Value = 532.375
$ws.Range("J127").Value = 1287.25
$ws.Range("K127").Value = 1597.125
$ws.Range("L127").Value = 3861.75
$ws.Range("M127").Value = 3362.875
$ws.Range("N127").Value = -13781.75

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 1025.6786
$ws.Range("J129").Value = 1089.3846
$ws.Range("L129").Value = 3268.1538
$ws.Range("N129").Value = -13268.1538

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 733185.7
$ws.Range("I135").Value = 380.65384
$ws.Range("J135").Value = 2003381
$ws.Range("K135").Value = 3425.88456
$ws.Range("L135").Value = 18030429
$ws.Range("M135").Value = -890.88456
$ws.Range("N135").Value = -18035499

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 28573224
$ws.Range("I137").Value = 1214.0625
$ws.Range("J137").Value = 333341340
$ws.Range("K137").Value = 3642.1875
$ws.Range("L137").Value = 1000024020
$ws.Range("M137").Value = -1092.1875
$ws.Range("N137").Value = -1000029120

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 10005380
$ws.Range("I138").Value = 12505475
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 37516425
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = -37511285
$ws.Range("N138").Value = -25280

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1114.2333
$ws.Range("I141").Value = 494.65216
$ws.Range("J141").Value = 3150
$ws.Range("K141").Value = 1483.95648
$ws.Range("L141").Value = 9450
$ws.Range("M141").Value = 3696.04352
$ws.Range("N141").Value = -19810

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1628.48
$ws.Range("I61").Value = 1585.6
$ws.Range("J61").Value = 1800
$ws.Range("K61").Value = 1585.6
$ws.Range("L61").Value = 1800
$ws.Range("M61").Value = -1373.6
$ws.Range("N61").Value = -2224

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 7467.409
$ws.Range("I74").Value = 1498.5
$ws.Range("J74").Value = 12441.5
$ws.Range("K74").Value = 1498.5
$ws.Range("L74").Value = 12441.5
$ws.Range("M74").Value = -624.5
$ws.Range("N74").Value = -14189.5

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 7467.409
$ws.Range("I77").Value = 1498.5
$ws.Range("J77").Value = 12441.5
$ws.Range("K77").Value = 7492.5
$ws.Range("L77").Value = 62207.5
$ws.Range("M77").Value = -3124.5
$ws.Range("N77").Value = -70943.5

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1628.48
$ws.Range("I136").Value = 1585.6
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 4756.799999999999
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = -2206.799999999999
$ws.Range("N136").Value = -10500

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 2032.0526
$ws.Range("I105").Value = 2159.2144
$ws.Range("J105").Value = 1676
$ws.Range("K105").Value = 2159.2144
$ws.Range("L105").Value = 1676
$ws.Range("M105").Value = -412.2143999999998
$ws.Range("N105").Value = -5170

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1255.1957
$ws.Range("I31").Value = 925
$ws.Range("J31").Value = 2093.3845
$ws.Range("K31").Value = 925
$ws.Range("L31").Value = 2093.3845
$ws.Range("M31").Value = -630
$ws.Range("N31").Value = -2683.3845

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1255.1957
$ws.Range("I34").Value = 925
$ws.Range("J34").Value = 2093.3845
$ws.Range("K34").Value = 925
$ws.Range("L34").Value = 2093.3845
$ws.Range("M34").Value = -723
$ws.Range("N34").Value = -2497.3845

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2389.9575
$ws.Range("I58").Value = 883.6
$ws.Range("J58").Value = 4101.727
$ws.Range("K58").Value = 883.6
$ws.Range("L58").Value = 4101.727
$ws.Range("M58").Value = -680.6
$ws.Range("N58").Value = -4507.727

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 2099.1667
$ws.Range("I99").Value = 1884.2858
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 1884.2858
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -386.2858000000001
$ws.Range("N99").Value = -5396

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 2099.1667
$ws.Range("I126").Value = 1884.2858
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5652.857400000001
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -3182.857400000001
$ws.Range("N126").Value = -12140

# Row 131: An Integral Reward / Integral Necklace of Crafting
$ws.Range("H131").Value = 21000
$ws.Range("J131").Value = 21000
$ws.Range("L131").Value = 21000
$ws.Range("N131").Value = -31080

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 1578.098
$ws.Range("I132").Value = 1491.9048
$ws.Range("J132").Value = 1980.3334
$ws.Range("K132").Value = 4475.7144
$ws.Range("L132").Value = 5941.0002
$ws.Range("M132").Value = -1945.7144
$ws.Range("N132").Value = -11001.0002

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 4157.5386
$ws.Range("I134").Value = 4713.4546
$ws.Range("J134").Value = 1100
$ws.Range("K134").Value = 14140.3638
$ws.Range("L134").Value = 3300
$ws.Range("M134").Value = -11605.3638
$ws.Range("N134").Value = -8370

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2389.9575
$ws.Range("I136").Value = 883.6
$ws.Range("J136").Value = 4101.727
$ws.Range("K136").Value = 2650.8
$ws.Range("L136").Value = 12305.181
$ws.Range("M136").Value = -100.8000000000002
$ws.Range("N136").Value = -17405.181

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 15437887
$ws.Range("I4").Value = 15437887
$ws.Range("K4").Value = 46313661
$ws.Range("M4").Value = -46313549

# Row 35: Whirled Peas / Pea Soup
$ws.Range("H35").Value = 3833.3333
$ws.Range("J35").Value = 3833.3333
$ws.Range("L35").Value = 11499.9999
$ws.Range("N35").Value = -12075.9999

# Row 57: The Egg Files / Deviled Eggs
$ws.Range("H57").Value = 1600
$ws.Range("J57").Value = 3000
$ws.Range("L57").Value = 9000
$ws.Range("N57").Value = -10118

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 3894.0303
$ws.Range("I131").Value = 6645.125
$ws.Range("K131").Value = 19935.375
$ws.Range("M131").Value = -14895.375

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 55556800
$ws.Range("I132").Value = 76924040
$ws.Range("J132").Value = 1981
$ws.Range("K132").Value = 692316360
$ws.Range("L132").Value = 17829
$ws.Range("M132").Value = -692313830
$ws.Range("N132").Value = -22889

# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 5311.125
$ws.Range("I134").Value = 5311.125
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15933.375
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -10863.375
$ws.Range("N134").ClearContents()

# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 48260.38
$ws.Range("I137").Value = 2561.4614
$ws.Range("J137").Value = 66825.56
$ws.Range("K137").Value = 7684.3842
$ws.Range("L137").Value = 200476.68
$ws.Range("M137").Value = -2584.3842
$ws.Range("N137").Value = -210676.68

$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 1963.3334
$ws.Range("I82").Value = 1450
$ws.Range("J82").Value = 2220
$ws.Range("K82").Value = 1450
$ws.Range("L82").Value = 2220
$ws.Range("M82").Value = -1089
$ws.Range("N82").Value = -2942

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 1963.3334
$ws.Range("I85").Value = 1450
$ws.Range("J85").Value = 2220
$ws.Range("K85").Value = 1450
$ws.Range("L85").Value = 2220
$ws.Range("M85").Value = -202
$ws.Range("N85").Value = -4716

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 1708.2632
$ws.Range("I136").Value = 894
$ws.Range("K136").Value = 2682
$ws.Range("M136").Value = -132

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 5301.864
$ws.Range("I81").Value = 1425.25
$ws.Range("J81").Value = 6163.3335
$ws.Range("K81").Value = 2850.5
$ws.Range("L81").Value = 12326.667
$ws.Range("M81").Value = -1789.5
$ws.Range("N81").Value = -14448.667

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 5301.864
$ws.Range("I84").Value = 1425.25
$ws.Range("J84").Value = 6163.3335
$ws.Range("K84").Value = 14252.5
$ws.Range("L84").Value = 61633.335
$ws.Range("M84").Value = -8948.5
$ws.Range("N84").Value = -72241.33499999999

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 4050.9812
$ws.Range("I132").Value = 4530.048
$ws.Range("J132").Value = 2221.818
$ws.Range("K132").Value = 13590.144
$ws.Range("L132").Value = 6665.454000000001
$ws.Range("M132").Value = -11060.144
$ws.Range("N132").Value = -11725.454

# Row 140: Glamorous Gloves / Thunderyards Silk Gloves of Casting
$ws.Range("H140").Value = 17943
$ws.Range("J140").Value = 17943
$ws.Range("L140").Value = 17943
$ws.Range("N140").Value = -28303
